$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 838, shifting existing rows 838:879 down to 839:880
$ws.Rows.Item(838).Insert()

# Populate the newly inserted row 838 with the new data point.
# Column A holds dates stored as plain text, so force a text number
# format before assigning the value (otherwise Excel auto-converts the
# "yyyy/mm/dd"-looking string into a date serial number). Reset the
# style back to Normal afterwards so the cell matches its neighbours,
# which carry no explicit style / number format of their own.
$ws.Range("A838").NumberFormat = "@"
$ws.Range("A838").Value = "2026/02/18"
$ws.Range("A838").Style = "Normal"
$ws.Range("B838").Value = "水"
$ws.Range("C838").Value = 18
$ws.Range("D838").Value = 201
